$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26: Jira id column changes from "N" to "Y" ---
$ws.Cells.Item(26,4).Value = "Y"

# --- Insert/populate new rows 27-29, copying formatting from row 26 first ---
$ws.Range("A26:E26").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A29:E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 27
$ws.Cells.Item(27,1).Value = "ENWIAM42"
$ws.Cells.Item(27,2).Value = "OPQA-3689|| OPQA-2172"
$ws.Cells.Item(27,3).Value = "Verify that the User is able to see 'Did you know? ...' Modal is displayed when user navigates from neon to ENW if Neon user has email same as existing steam acount (used to login ENW).`nVerify that User should be taken back to Neon Home page, When User clicks on Close '[X] ' button on linking modals while Navigation from Neon via clicking link to EndNote"
$ws.Cells.Item(27,4).Value = "Y"
$ws.Rows.Item(27).RowHeight = 90

# Row 28
$ws.Cells.Item(28,1).Value = "ENWIAM50"
$ws.Cells.Item(28,2).Value = "OPQA-1707||OPQA-1734"
$ws.Cells.Item(28,3).Value = "Verify that Neon Landing page, displays Neon branding , marketing copy , New icon and also integration with Endnote"
$ws.Cells.Item(28,4).Value = "Y"
$ws.Rows.Item(28).RowHeight = 30

# Row 29
$ws.Cells.Item(29,1).Value = "ENWIAM51"
$ws.Cells.Item(29,2).Value = "OPQA-1673 || OPQA-1681 || OPQA-1691 || OPQA-1817 || OPQA-3648 || OPQA-3649"
$ws.Cells.Item(29,3).Value = "Verify that on ENW landing page displays,EndNote branding and marketing copy and integration with Project Neon"
$ws.Cells.Item(29,4).Value = "Y"
$ws.Rows.Item(29).RowHeight = 60

# Rows 28 and 29 use wrap-text style for column B (matches style s="7" in target)
$ws.Cells.Item(28,2).WrapText = $true
$ws.Cells.Item(29,2).WrapText = $true

# Update selection / active cell to match final view state
$ws.Activate()
$ws.Range("C27").Select()

